# 22 scenarios done- login page
#
# The login-data sheet's B2 cell held a hyperlink display-text that duplicated
# the "herokuapp1.com" typo URL already used elsewhere. Retype it with the
# corrected display text and apply the built-in "Hyperlink" cell style (left/top,
# wrapped), then leave the selection on B2 - matching what a user does when they
# click into B2, correct the URL text, and restyle it from the Cell Styles gallery.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B2")

# New display text for the cell (the underlying hyperlink target is untouched).
$cell.Value = "https://llms-frontend-api-hackathon-apr-326235f3973d.herokuapp.com/"

# Apply the workbook's built-in Hyperlink cell style, then re-assert the
# left/top wrapped alignment the rest of the row already uses.
$cell.Style = "Hyperlink"
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160
$cell.WrapText = $true

# Leave the selection on B2.
$null = $cell.Select()
